# Natmi following Dr Hou advice
# Rebuild the Cd200 -> Cd200r4 LR-pair table with the full sending/target
# cluster cross-join (ECs, FAPs, M2, sCs) per the updated analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd200"
$ws.Range("C2").Value = "Cd200r4"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 65.569613
$ws.Range("H2").Value = 196.708839
$ws.Range("I2").Value = 0.6815983651189281
$ws.Range("J2").Value = 0.681598365118928
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1322443333333333
$ws.Range("N2").Value = 0.396733
$ws.Range("O2").Value = 0.01281473224718535
$ws.Range("P2").Value = 0.01281473224718535
$ws.Range("Q2").Value = 8.671209758109667
$ws.Range("R2").Value = 78.040887822987
$ws.Range("S2").Value = 0.00873450054911834
$ws.Range("T2").Value = 0.008734500549118339

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd200"
$ws.Range("C3").Value = "Cd200r4"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 65.569613
$ws.Range("H3").Value = 196.708839
$ws.Range("I3").Value = 0.6815983651189281
$ws.Range("J3").Value = 0.681598365118928
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 10.18746666666667
$ws.Range("N3").Value = 30.5624
$ws.Range("O3").Value = 0.9871852677528148
$ws.Range("P3").Value = 0.9871852677528146
$ws.Range("Q3").Value = 667.9882467837334
$ws.Range("R3").Value = 6011.894221053601
$ws.Range("S3").Value = 0.6728638645698098
$ws.Range("T3").Value = 0.6728638645698096

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cd200"
$ws.Range("C4").Value = "Cd200r4"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 4.561623333333333
$ws.Range("H4").Value = 13.68487
$ws.Range("I4").Value = 0.04741823024467683
$ws.Range("J4").Value = 0.04741823024467683
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1322443333333333
$ws.Range("N4").Value = 0.396733
$ws.Range("O4").Value = 0.01281473224718535
$ws.Range("P4").Value = 0.01281473224718535
$ws.Range("Q4").Value = 0.6032488366344444
$ws.Range("R4").Value = 5.42923952971
$ws.Range("S4").Value = 0.0006076519242209196
$ws.Range("T4").Value = 0.0006076519242209196

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd200"
$ws.Range("C5").Value = "Cd200r4"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 4.561623333333333
$ws.Range("H5").Value = 13.68487
$ws.Range("I5").Value = 0.04741823024467683
$ws.Range("J5").Value = 0.04741823024467683
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 10.18746666666667
$ws.Range("N5").Value = 30.5624
$ws.Range("O5").Value = 0.9871852677528148
$ws.Range("P5").Value = 0.9871852677528146
$ws.Range("Q5").Value = 46.47138565422222
$ws.Range("R5").Value = 418.242470888
$ws.Range("S5").Value = 0.04681057832045591
$ws.Range("T5").Value = 0.0468105783204559

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cd200"
$ws.Range("C6").Value = "Cd200r4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2.0
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.611825
$ws.Range("H6").Value = 1.835475
$ws.Range("I6").Value = 0.006359941757455365
$ws.Range("J6").Value = 0.006359941757455365
$ws.Range("K6").Value = 1.0
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1322443333333333
$ws.Range("N6").Value = 0.396733
$ws.Range("O6").Value = 0.01281473224718535
$ws.Range("P6").Value = 0.01281473224718535
$ws.Range("Q6").Value = 0.08091038924166666
$ws.Range("R6").Value = 0.728193503175
$ws.Range("S6").Value = 0.00008150095072948391
$ws.Range("T6").Value = 0.00008150095072948391

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cd200"
$ws.Range("C7").Value = "Cd200r4"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 2.0
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.611825
$ws.Range("H7").Value = 1.835475
$ws.Range("I7").Value = 0.006359941757455365
$ws.Range("J7").Value = 0.006359941757455365
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 10.18746666666667
$ws.Range("N7").Value = 30.5624
$ws.Range("O7").Value = 0.9871852677528148
$ws.Range("P7").Value = 0.9871852677528146
$ws.Range("Q7").Value = 6.232946793333333
$ws.Range("R7").Value = 56.09652114
$ws.Range("S7").Value = 0.006278440806725882
$ws.Range("T7").Value = 0.00627844080672588

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cd200"
$ws.Range("C8").Value = "Cd200r4"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 25.456719
$ws.Range("H8").Value = 76.37015699999999
$ws.Range("I8").Value = 0.2646234628789398
$ws.Range("J8").Value = 0.2646234628789398
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1322443333333333
$ws.Range("N8").Value = 0.396733
$ws.Range("O8").Value = 0.01281473224718535
$ws.Range("P8").Value = 0.01281473224718535
$ws.Range("Q8").Value = 3.366506833008999
$ws.Range("R8").Value = 30.298561497081
$ws.Range("S8").Value = 0.003391078823116604
$ws.Range("T8").Value = 0.003391078823116604

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cd200"
$ws.Range("C9").Value = "Cd200r4"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 25.456719
$ws.Range("H9").Value = 76.37015699999999
$ws.Range("I9").Value = 0.2646234628789398
$ws.Range("J9").Value = 0.2646234628789398
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 10.18746666666667
$ws.Range("N9").Value = 30.5624
$ws.Range("O9").Value = 0.9871852677528148
$ws.Range("P9").Value = 0.9871852677528146
$ws.Range("Q9").Value = 259.3394762552
$ws.Range("R9").Value = 2334.0552862968
$ws.Range("S9").Value = 0.2612323840558233
$ws.Range("T9").Value = 0.2612323840558232

